# Update cryptocurrency Price column (D) with refreshed values.
# Source data was scraped/updated on Wed Dec 14 19:34:16 UTC 2022 (GitHub Actions).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column D holds values as plain text (e.g. "270.90") so that
# trailing zeros / exact scraped formatting survive. Force each target
# cell to Text format before writing so Excel doesn't silently coerce
# the numeric-looking string into a floating point number.
$priceUpdates = @{
    "D2" = "271.25"
    "D3" = "22.94"
    "D4" = "6.376"
    "D5" = "0.06220"
    "D6" = "3.650"
    "D7" = "6.763"
    "D8" = "1.405"
    "D9" = "0.8374"
    "D10" = "0.01366"
    "D12" = "0.08308"
    "D13" = "0.03404"
    "D14" = "0.03184"
    "D15" = "0.09295"
    "D16" = "3.926"
    "D17" = "0.001729"
    "D18" = "0.04868"
    "D19" = "0.006234"
    "D20" = "0.005422"
    "D21" = "0.001096"
    "D22" = "0.0001508"
    "D23" = "3.753"
    "D24" = "2.333"
    "D25" = "0.3339"
    "D26" = "0.1252"
    "D27" = "0.0002697"
    "D40" = "0.04678"
    "D41" = "0.006931"
    "D42" = "0.1164"
    "D43" = "0.003479"
    "D44" = "0.01235"
    "D45" = "0.00006303"
    "D46" = "0.00000000754"
    "D47" = "0.7038"
    "D48" = "0.1300"
    "D49" = "0.00002111"
    "D50" = "0.01247"
}

foreach ($cellRef in $priceUpdates.Keys) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $priceUpdates[$cellRef]
}
